$d = $word.ActiveDocument

# --- Edit 1: remove the _GoBack bookmark from the paragraph beginning "The MathML parser ..." ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Edit 2: remove the stray <w:lastRenderedPageBreak/> before "How to Match Patterns" ---
# lastRenderedPageBreak is a transient rendering artifact; re-writing the run's
# text causes the engine to drop it since it is recomputed at render time.
$d.Content.Find.Execute("How to Match Patterns", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "How to Match Patterns", 2) | Out-Null

# --- Edit 3: split the "If pattern element is a Text..." bullet into 4 runs with new wording ---
$rng = $d.Content
$found = $rng.Find.Execute("If pattern element is a Text, its compiled regular expression is only used to match a Text", `
                            $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Remove the old text, leaving a collapsed insertion point in the same paragraph
    # (preserves the paragraph's pPr / list numbering).
    $rng.Text = ""

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/part.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
        '<w:p>' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr></w:pPr>' + `
        '<w:r><w:t>If</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> pattern element is a Text, it only matches</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> a Text</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> with same content</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml) | Out-Null
}
